# "sec key working at level 1"
#
# Two logical changes, matching the OOXML diff:
#
#   1) The provenance sheet's "timestamp" row (B12) is re-stamped with a
#      later build time (the workbook was regenerated).
#   2) The "08-BC" sheet's security-key cell (A3) is advanced from the
#      base key ("08-BC") through the level-1/level-2/level-3 key labels,
#      ending up on "08-BC.S-03" - the shared string table picks up the
#      "08-BC.S-01" / "08-BC.S-02" / "08-BC.S-03" labels as the key is
#      advanced, and the cell ends up referencing the final label.

$wb = $excel.ActiveWorkbook

# --- 1) Re-stamp the provenance timestamp ---
$wsProv = $wb.Worksheets.Item("provenance")
$wsProv.Range("B12").Value = 43435.42160976339

# --- 2) Advance the "08-BC" sheet's security key to level 3 ---
$wsBC = $wb.Worksheets.Item("08-BC")
$wsBC.Range("A3").Value = "08-BC.S-01"
$wsBC.Range("A3").Value = "08-BC.S-02"
$wsBC.Range("A3").Value = "08-BC.S-03"
